# Aufwand_git.xlsx update
# - Add two new time-tracking entries (rows 60 & 61) on sheet "Tabelle1"
# - Shared strings / formulas / dimension / view selection update automatically

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# New row 60: 2024-03-19 (serial 45370), 3 hours, "Presentation Vorb"
$ws.Cells.Item(60, 1).Value = 45370
$ws.Cells.Item(60, 1).NumberFormat = "d-mmm"
$ws.Cells.Item(60, 2).Value = 3
$ws.Cells.Item(60, 3).Value = "Presentation Vorb"

# New row 61: 2024-03-19 (serial 45370), 3 hours, "Remove Supp and Conf, ShaclOrItems, Default Shapes saved"
$ws.Cells.Item(61, 1).Value = 45370
$ws.Cells.Item(61, 1).NumberFormat = "d-mmm"
$ws.Cells.Item(61, 2).Value = 3
$ws.Cells.Item(61, 3).Value = "Remove Supp and Conf, ShaclOrItems, Default Shapes saved"

# Scroll the view so the newly added rows are visible, and move the
# selection to match the author's final cursor position.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 56
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F64").Select()
